$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B9").Value = 3831655.92
$ws.Range("C9").Value = 604610.24
$ws.Range("D9").Value = 4436266.16
$ws.Range("E9").Value = 13.62880896217462
$ws.Range("F9").Value = 86.37119103782537
$ws.Range("G9").Value = -41.5674350296405
$ws.Range("H9").Value = -30.80559035140188
$ws.Range("I9").Value = 38736
$ws.Range("J9").Value = 1657
$ws.Range("K9").Value = 40393
$ws.Range("L9").Value = 27942
$ws.Range("M9").Value = 158.7669515424809
$ws.Range("N9").Value = 8.39343154035992
